$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll + selection change ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("F21").Select()

# --- New helper cell E20 = 0 (prefixes the generated id with a leading 0) ---
$ws.Range("E20").Value = 0

# --- Existing formulas now also reference $E$20, re-grouped as shared formulas ---
$ws.Range("F1:F3").Formula = '=CONCATENATE($J$1,$E$20,A1,$J$3,B1,$J$3,C1,$J$2)'
$ws.Range("F7:F8").Formula = '=CONCATENATE($J$1,$E$20,A7,$J$3,B7,$J$3,C7,$J$2)'
$ws.Range("F12").Formula = '=CONCATENATE($J$1,$E$20,A12,$J$3,B12,$J$3,C12,$J$2)'

# --- New data rows 21-43 ---
$data = @(
    @(21, 1238471, 1203, 7),
    @(22, 1238471, 8892, 5),
    @(23, 2321423, 2056, 8),
    @(24, 2321423, 2984, 8),
    @(25, 2321423, 1784, 7),
    @(26, 4528483, 1203, 10),
    @(27, 4528483, 8892, 7),
    @(28, 4737492, 3282, 8),
    @(29, 4737492, 6234, 5),
    @(30, 4917493, 1203, 6),
    @(31, 9373493, 8892, 8),
    @(32, 9373493, 5738, 8),
    @(33, 9373493, 7231, 9),
    @(34, 5838743, 3627, 10),
    @(35, 4782107, 3492, 6),
    @(36, 9043278, 7329, 5),
    @(37, 4309871, 6234, 7),
    @(38, 9043278, 6321, 8),
    @(39, 9373493, 3229, 9),
    @(40, 9043278, 7482, 10),
    @(41, 2374827, 7482, 10),
    @(42, 3283295, 6321, 9),
    @(43, 3421987, 6321, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
}

# F21 starts as a standalone formula, F22:F43 become one shared group
$ws.Range("F21").Formula = '=CONCATENATE($J$1,$E$20,A21,$J$3,B21,$J$3,C21,$J$2)'
$ws.Range("F22:F43").Formula = '=CONCATENATE($J$1,$E$20,A22,$J$3,B22,$J$3,C22,$J$2)'
